$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: column A holds sequential index values (continuing 204..215),
# column B holds the corresponding normalized values, appended to rows 206-217.
$aValues = @(204, 205, 206, 207, 208, 209, 210, 211, 212, 213, 214, 215)
$bValues = @(
    [double]"-1.572815951552305E-16",
    [double]"-1.160687707562664E-16",
    [double]"-2.109423746787797E-16",
    [double]"-9.251858538542972E-17",
    [double]"1.040834085586084E-17",
    [double]"-1.586032892321652E-16",
    [double]"1.480297366166875E-16",
    [double]"3.885780586188048E-17",
    [double]"-3.469446951953614E-18",
    [double]"0",
    [double]"-8.326672684688674E-17",
    [double]"0"
)

$styleSourceCell = $ws.Cells.Item(205, 1)
$startRow = 206
for ($i = 0; $i -lt $aValues.Count; $i++) {
    $row = $startRow + $i

    # Copy the formatting (bold/border/centered style) from the last existing
    # "index" cell in column A onto the newly appended one.
    $styleSourceCell.Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $aValues[$i]
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}

$excel.CutCopyMode = $false
